# Apply data corrections to the CI-config parsing table ("arreglos en
# 'parseo' de ficheros YML de CI") and drop the trailing rows that no
# longer have any matches (Objective-C, Vim script, Jsonnet, Dart, PHP,
# SCSS), shrinking the used range from A1:N30 down to A1:N24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - None
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3

# Row 4 - TypeScript
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 0

# Row 5 - Python
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 0

# Row 6 - JavaScript
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 30

# Row 7 - C++
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 5
$ws.Range("I7").Value = 0

# Row 8 - Rust
$ws.Range("E8").Value = 4

# Row 9 - Vue
$ws.Range("C9").Value = 3
$ws.Range("E9").Value = 3

# Row 10 - C#
$ws.Range("E10").Value = 2

# Row 11 - Shell
$ws.Range("E11").Value = 4

# Row 12 - Java
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 9

# Row 13 - CSS
$ws.Range("E13").Value = 0

# Row 14 - Kotlin
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 2

# Row 15 - Go
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 8
$ws.Range("G15").Value = 0
$ws.Range("L15").Value = 0

# Row 16 - Ruby
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 2

# Row 17 - Swift
$ws.Range("E17").Value = 5

# Row 18 - C
$ws.Range("E18").Value = 1

# Row 19 - Jupyter Notebook
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2

# Row 23 - HTML
$ws.Range("E23").Value = 2

# Remove trailing rows 25-30 (Objective-C, Vim script, Jsonnet, Dart,
# PHP, SCSS) - they no longer have any entries worth reporting.
$ws.Range("A25:N30").EntireRow.Delete() | Out-Null
